# Added Panel Accessories Test Data For Spain/Turkey/Hungary market
#
# Clones the "Greece" sheet (the template row/column/style layout used by
# every market tab in this workbook) three times to create "Hungary",
# "Spain" and "Turkey" tabs, fills in their market name / ticket-id cells,
# and updates the view/selection state to match what Excel leaves behind
# after this kind of edit (Turkey ends up the active/selected tab).

$wb = $excel.ActiveWorkbook
$greece = $wb.Worksheets.Item("Greece")

# --- Hungary -------------------------------------------------------------
$greece.Copy($null, $greece)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-3104/T2998/T2980/T2994"

# --- Spain -----------------------------------------------------------------
$greece.Copy($null, $hungary)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2037/T2051/T2056"

# --- Turkey ----------------------------------------------------------------
$greece.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3311/T3317/T3300"

# --- View / selection state ------------------------------------------------
# Greece is no longer the tab shown/selected; the whole data range is
# selected instead of the old B4 pick, scrolled down a bit.
$greece.Activate()
$greece.Range("A1:D22").Select()

$hungary.Activate()
$hungary.Range("A1:D22").Select()

$spain.Activate()
$spain.Range("A1:D22").Select()

# Turkey ends up the active/visible tab, with cell I4 selected.
$turkey.Activate()
$turkey.Range("I4").Select()
